# DTFS2-8124: accrual currency validation changed to ISO4217
#
# Adds two new columns (L: "Accrual currency", M: "accrual exchange rate")
# to the fixture's first worksheet, with sample/invalid data in the two
# existing data rows so the fixture can be used to test ISO4217 validation
# of the new "accrual currency" field.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 - new headers for columns L and M
$ws.Range("L1").Value = "Accrual currency"
$ws.Range("M1").Value = "accrual exchange rate"

# Row 2 - sample data: an invalid (non-ISO4217) accrual currency code
$ws.Range("L2").Value = "INRA"
$ws.Range("M2").Value = 1.223

# Row 3 - sample data: another invalid (non-ISO4217) accrual currency code
$ws.Range("L3").Value = "A"
$ws.Range("M3").Value = 2.33

# Reflect the view state captured when the fixture was last edited: the
# new M column is selected and the sheet is scrolled right so column B
# is the first visible column.
$ws.Range("M4").Select()
$excel.ActiveWindow.ScrollColumn = 2
